$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 158.36363
$ws.Range("I6").Value = 154.2
$ws.Range("K6").Value = 462.6
$ws.Range("M6").Value = -350.6

$ws.Range("H11").Value = 317
$ws.Range("I11").Value = 317
$ws.Range("K11").Value = 317
$ws.Range("M11").Value = -177

$ws.Range("H33").Value = 136.38461
$ws.Range("I33").Value = 134.91667
$ws.Range("K33").Value = 134.91667
$ws.Range("M33").Value = 94.08332999999999

$ws.Range("H100").Value = 2493.125
$ws.Range("I100").Value = 2157.8333
$ws.Range("K100").Value = 2157.8333
$ws.Range("M100").Value = -1616.8333

$ws.Range("H116").Value = 6599.6
$ws.Range("I116").Value = 4332.6665
$ws.Range("K116").Value = 4332.6665
$ws.Range("M116").Value = -890.6665000000003

$ws.Range("H132").Value = 1125.238
$ws.Range("I132").Value = 1179.9474
$ws.Range("K132").Value = 3539.8422
$ws.Range("M132").Value = -1009.8422

$ws.Range("H135").Value = 1114
$ws.Range("I135").Value = 1002.0909
$ws.Range("K135").Value = 9018.8181
$ws.Range("M135").Value = -6483.8181

$ws.Range("H137").Value = 1730.1765
$ws.Range("I137").Value = 1737.7142
$ws.Range("K137").Value = 5213.142599999999
$ws.Range("M137").Value = -2663.142599999999

$ws.Range("H141").Value = 1684.6072
$ws.Range("I141").Value = 1623.5555
$ws.Range("J141").Value = 3333
$ws.Range("K141").Value = 4870.666499999999
$ws.Range("L141").Value = 9999
$ws.Range("M141").Value = 309.3335000000006
$ws.Range("N141").Value = -20359

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1610.5555
$ws.Range("I2").Value = 1319.2
$ws.Range("K2").Value = 1319.2
$ws.Range("M2").Value = -1206.2

$ws.Range("H32").Value = 4182.923
$ws.Range("I32").Value = 2554.8044
$ws.Range("J32").Value = 16665.166
$ws.Range("K32").Value = 2554.8044
$ws.Range("L32").Value = 16665.166
$ws.Range("M32").Value = -2267.8044
$ws.Range("N32").Value = -17239.166

$ws.Range("H45").Value = 29973.5
$ws.Range("I45").Value = 819.7143
$ws.Range("K45").Value = 819.7143
$ws.Range("M45").Value = -442.7143

$ws.Range("H110").Value = 20002
$ws.Range("I110").Value = 20002
$ws.Range("K110").Value = 20002
$ws.Range("M110").Value = -17957

$ws.Range("H116").Value = 1610.5555
$ws.Range("I116").Value = 1319.2
$ws.Range("K116").Value = 1319.2
$ws.Range("M116").Value = 974.8

$ws.Range("H122").Value = 5624.45
$ws.Range("I122").Value = 4922.231
$ws.Range("J122").Value = 6928.5713
$ws.Range("K122").Value = 14766.693
$ws.Range("L122").Value = 20785.7139
$ws.Range("M122").Value = -12316.693
$ws.Range("N122").Value = -25685.7139

$ws.Range("H132").Value = 1012.41174
$ws.Range("I132").Value = 1012.41174
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3037.23522
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -507.23522
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1610.5555
$ws.Range("I3").Value = 1319.2
$ws.Range("K3").Value = 1319.2
$ws.Range("M3").Value = -1205.2

$ws.Range("H134").Value = 1794.8
$ws.Range("I134").Value = 1794.8
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5384.4
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2849.4
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 50000976
$ws.Range("I2").Value = 50000976
$ws.Range("K2").Value = 50000976
$ws.Range("M2").Value = -50000863

$ws.Range("H22").Value = 334
$ws.Range("I22").Value = 313.5
$ws.Range("J22").Value = 375
$ws.Range("K22").Value = 313.5
$ws.Range("L22").Value = 375
$ws.Range("M22").Value = 36.5
$ws.Range("N22").Value = -1075

$ws.Range("H86").Value = 8333.333000000001
$ws.Range("J86").Value = 16000
$ws.Range("L86").Value = 16000
$ws.Range("N86").Value = -18246

$ws.Range("H89").Value = 8333.333000000001
$ws.Range("J89").Value = 16000
$ws.Range("L89").Value = 80000
$ws.Range("N89").Value = -91232

$ws.Range("H98").Value = 125000
$ws.Range("J98").Value = 125000
$ws.Range("L98").Value = 125000
$ws.Range("N98").Value = -129492

$ws.Range("H131").Value = 80000
$ws.Range("J131").Value = 80000
$ws.Range("L131").Value = 80000
$ws.Range("N131").Value = -90080

$ws.Range("H132").Value = 2247.6
$ws.Range("I132").Value = 1947
$ws.Range("J132").Value = 3450
$ws.Range("K132").Value = 5841
$ws.Range("L132").Value = 10350
$ws.Range("M132").Value = -3311
$ws.Range("N132").Value = -15410

$ws.Range("H134").Value = 3988
$ws.Range("I134").Value = 3980.6667
$ws.Range("J134").Value = 3999
$ws.Range("K134").Value = 11942.0001
$ws.Range("L134").Value = 11997
$ws.Range("M134").Value = -9407.000100000001
$ws.Range("N134").Value = -17067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 900
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 900
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 2700
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -3148

$ws.Range("H28").Value = 2580
$ws.Range("I28").Value = 2580
$ws.Range("K28").Value = 7740
$ws.Range("M28").Value = -7508

$ws.Range("H37").Value = 88900
$ws.Range("J37").Value = 88900
$ws.Range("L37").Value = 266700
$ws.Range("N37").Value = -266924

$ws.Range("H109").Value = 1000.8889
$ws.Range("I109").Value = 1001.125
$ws.Range("K109").Value = 3003.375
$ws.Range("M109").Value = -1963.375

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws.Range("H121").Value = 850.2
$ws.Range("J121").Value = 855.7778
$ws.Range("L121").Value = 2567.3334
$ws.Range("N121").Value = -5187.3334

$ws.Range("H129").Value = 1483.5
$ws.Range("I129").Value = 363.75
$ws.Range("K129").Value = 1091.25
$ws.Range("M129").Value = 3908.75

$ws.Range("H132").Value = 5103.647
$ws.Range("I132").Value = 5359.6924
$ws.Range("J132").Value = 4271.5
$ws.Range("K132").Value = 48237.2316
$ws.Range("L132").Value = 38443.5
$ws.Range("M132").Value = -45707.2316
$ws.Range("N132").Value = -43503.5

$ws.Range("H137").Value = 6000
$ws.Range("J137").Value = 5500
$ws.Range("L137").Value = 16500
$ws.Range("N137").Value = -26700

$ws.Range("H140").Value = 2697
$ws.Range("I140").Value = 1929.3334
$ws.Range("K140").Value = 5788.0002
$ws.Range("M140").Value = -608.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 62655.938
$ws.Range("I3").Value = 62655.938
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 62655.938
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -62539.938
$ws.Range("N3").ClearContents()

$ws.Range("H132").Value = 3282.25
$ws.Range("I132").Value = 3043
$ws.Range("K132").Value = 9129
$ws.Range("M132").Value = -6599

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1877.5
$ws.Range("I7").Value = 1750
$ws.Range("K7").Value = 1750
$ws.Range("M7").Value = -1638

$ws.Range("H43").Value = 2039490.4
$ws.Range("J43").Value = 2039490.4
$ws.Range("L43").Value = 2039490.4
$ws.Range("N43").Value = -2039876.4

$ws.Range("H82").Value = 3034
$ws.Range("J82").Value = 749
$ws.Range("L82").Value = 749
$ws.Range("N82").Value = -1471

$ws.Range("H85").Value = 3034
$ws.Range("J85").Value = 749
$ws.Range("L85").Value = 749
$ws.Range("N85").Value = -3245

$ws.Range("H126").Value = 1877.5
$ws.Range("I126").Value = 1750
$ws.Range("K126").Value = 5250
$ws.Range("M126").Value = -2780

$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470

$ws.Range("H136").Value = 3299
$ws.Range("I136").Value = 3299
$ws.Range("K136").Value = 9897
$ws.Range("M136").Value = -7347

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 4005
$ws.Range("I7").Value = 4005
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4005
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3892
$ws.Range("N7").ClearContents()

$ws.Range("H29").Value = 25149.334
$ws.Range("I29").Value = 25149.334
$ws.Range("K29").Value = 25149.334
$ws.Range("M29").Value = -24859.334

$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H122").Value = 4999.6665
$ws.Range("I122").Value = 5499.8
$ws.Range("J122").Value = 2499
$ws.Range("K122").Value = 16499.4
$ws.Range("L122").Value = 7497
$ws.Range("M122").Value = -14049.4
$ws.Range("N122").Value = -12397

$ws.Range("H126").Value = 3093.7368
$ws.Range("I126").Value = 2901.5
$ws.Range("J126").Value = 3307.3333
$ws.Range("K126").Value = 8704.5
$ws.Range("L126").Value = 9921.999899999999
$ws.Range("M126").Value = -6234.5
$ws.Range("N126").Value = -14861.9999
